$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) cells we are about to overwrite stay as text values
# (these are number-like strings such as "0.9999" that Excel would
# otherwise auto-convert to a numeric value when assigned via .Value).
$ws.Range("D2:D20").NumberFormat = "@"
$ws.Range("D22:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.414.72'
$ws.Range("E2").Value = '  -0.69%  '

$ws.Range("D3").Value = '1.790.26'
$ws.Range("E3").Value = '  -1.92%  '

$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").Value = '340.93'
$ws.Range("E5").Value = '  +0.64%  '

$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").Value = '0.3953'
$ws.Range("E7").Value = '  +3.19%  '

$ws.Range("D8").Value = '0.3475'
$ws.Range("E8").Value = '  -1.70%  '

$ws.Range("D9").Value = '48.07'
$ws.Range("E9").Value = '  -3.39%  '

$ws.Range("D10").Value = '1.199'
$ws.Range("E10").Value = '  -3.38%  '

$ws.Range("D11").Value = '0.07502'
$ws.Range("E11").Value = '  -3.16%  '

$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = '21.86'
$ws.Range("E13").Value = '  -2.61%  '

$ws.Range("D14").Value = '6.497'
$ws.Range("E14").Value = '  -2.11%  '

$ws.Range("D15").Value = '1.781.55'
$ws.Range("E15").Value = '  -2.31%  '

$ws.Range("D16").Value = '7.127'
$ws.Range("E16").Value = '  -1.22%  '

$ws.Range("D17").Value = '0.00001098'
$ws.Range("E17").Value = '  -2.62%  '

$ws.Range("D18").Value = '0.06712'
$ws.Range("E18").Value = '  -0.13%  '

$ws.Range("D19").Value = '84.90'
$ws.Range("E19").Value = '  -2.72%  '

$ws.Range("D20").Value = '0.9985'
$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("E21").Value = '  +0.25%  '

$ws.Range("D22").Value = '6.521'
$ws.Range("E22").Value = '  -0.48%  '

$ws.Range("D23").Value = '27.346.67'
$ws.Range("E23").Value = '  -0.97%  '

$ws.Range("D24").Value = '12.43'
$ws.Range("E24").Value = '  -5.60%  '

$ws.Range("D25").Value = '2.387'
$ws.Range("E25").Value = '  -3.55%  '

$ws.Range("D26").Value = '21.31'
$ws.Range("E26").Value = '  -3.80%  '

$ws.Range("D27").Value = '1.464'
$ws.Range("E27").Value = '  -1.88%  '

$ws.Range("D28").Value = '2.502'
$ws.Range("E28").Value = '  -6.72%  '

$ws.Range("D29").Value = '157.79'
$ws.Range("E29").Value = '  +3.10%  '

$ws.Range("D30").Value = '1.985.75'
$ws.Range("E30").Value = '  -2.08%  '

$ws.Range("D31").Value = '136.19'
$ws.Range("E31").Value = '  +0.53%  '

$ws.Range("D32").Value = '4.027'
$ws.Range("E32").Value = '  -1.39%  '

$ws.Range("D33").Value = '6.009'
$ws.Range("E33").Value = '  -5.60%  '

$ws.Range("D34").Value = '0.08837'
$ws.Range("E34").Value = '  +0.18%  '

$ws.Range("D35").Value = '13.04'
$ws.Range("E35").Value = '  -6.75%  '

$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.618'
$ws.Range("E36").Value = '  -4.66%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02427'
$ws.Range("E37").Value = '  +0.66%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.06506'
$ws.Range("E38").Value = '  -0.41%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '5.425'
$ws.Range("E39").Value = '  -4.02%  '

$ws.Range("D40").Value = '0.6829'
$ws.Range("E40").Value = '  -3.27%  '

$ws.Range("D41").Value = '0.2217'
$ws.Range("E41").Value = '  -2.36%  '

$ws.Range("D42").Value = '1.252'
$ws.Range("E42").Value = '  -3.58%  '

$ws.Range("D43").Value = '8.382'
$ws.Range("E43").Value = '  -8.60%  '

$ws.Range("D44").Value = '14.49'
$ws.Range("E44").Value = '  -2.34%  '

$ws.Range("D45").Value = '0.9988'
$ws.Range("E45").Value = '  -0.09%  '

$ws.Range("D46").Value = '0.6398'
$ws.Range("E46").Value = '  -3.70%  '

$ws.Range("D47").Value = '3.880'
$ws.Range("E47").Value = '  -0.97%  '

$ws.Range("D48").Value = '2.143'
$ws.Range("E48").Value = '  -2.26%  '

$ws.Range("D49").Value = '132.36'
$ws.Range("E49").Value = '  -0.82%  '

$ws.Range("D50").Value = '0.07165'
$ws.Range("E50").Value = '  -1.94%  '

$ws.Range("D51").Value = '79.24'
$ws.Range("E51").Value = '  -3.13%  '
